# إضافة حدث جديد في Card23
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# --- Fill in the previously-blank "nan" placeholder cells on row 13 ---
# (columns B..K and M were empty inline strings; the generator backfills
#  them with the literal text "nan" once the row is no longer the last one)
$nanCols = @(2,3,4,5,6,7,8,9,10,11,13)  # B,C,D,E,F,G,H,I,J,K,M
foreach ($col in $nanCols) {
    $ws.Cells.Item(13, $col).Value = "nan"
}

# --- Add the new event as row 14 ---
# Make sure the row card id is stored as text, matching the rest of column A
$ws.Range("A14").NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "23"

# Columns B..K are intentionally left blank for the newly added event
$ws.Range("L14:O14").NumberFormat = "@"
$ws.Cells.Item(14, 12).Value = "10\8\2024"
$ws.Cells.Item(14, 13).Value = "4320 h"
$ws.Cells.Item(14, 14).Value = "تم تشحيم الكنه بالكامل+عمل صيانه"
$ws.Cells.Item(14, 15).Value = "تيم العمل"
